# Rapport nuevo con frases largas divididas en mensajes consecutivos
# Update the "DariJourney" sheet (R-series rapport-building scripts):
# - Rows 34/35 (R-12, R-11) get reworded B-column text.
# - Rows 36-45 (R-10 .. R-1) are replaced by a new longer sequence split
#   into consecutive "a/b" messages, pushing the tail down to row 50.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DariJourney")

# Create 5 new rows (46-50) with the same style as the existing R-series
# rows (s="8") by copying a same-sized styled block down.
$ws.Range("A41:C45").Copy($ws.Range("A46:C50"))

$ws.Range("A34").Value = 'R-12'
$ws.Range("B34").Value = 'que descanses lindo mor, portate bien ✨ mañana quiero que me cuentes que soñaste 💗'

$ws.Range("A35").Value = 'R-11'
$ws.Range("B35").Value = 'bueno papi, me encanta esta conversacion 🥰 aqui voy a estar cuando quieras hablar, si? me haces el dia ✨'

$ws.Range("A36").Value = 'R-10b'
$ws.Range("B36").Value = 'como que de verdad te importa lo que digo y eso se siente bonito 💗'
$ws.Range("C36").Value = 'Enviar inmediatamente despues de R-10a, sin esperar respuesta.'

$ws.Range("A37").Value = 'R-10a'
$ws.Range("B37").Value = 'es que la mayoria de la gente solo escribe por escribir, pero contigo es distinto...'
$ws.Range("C37").Value = 'Ego boost profundo. El es especial.'

$ws.Range("A38").Value = 'R-9'
$ws.Range("B38").Value = 'ay mor de verdad que hablar contigo me pone contenta ❤️ como que siento que me entiendes'
$ws.Range("C38").Value = 'Validacion emocional. ''Eres diferente'' implicito.'

$ws.Range("A39").Value = 'R-8b'
$ws.Range("B39").Value = 'yo prefiero una peliculita en la cama con comidita rica 💗'
$ws.Range("C39").Value = 'Enviar inmediatamente despues de R-8a, sin esperar respuesta.'

$ws.Range("A40").Value = 'R-8a'
$ws.Range("B40").Value = 'es que a mi me gustan las cosas tranquilitas, amor... nada de fiestas ni eso'
$ws.Range("C40").Value = 'Mostrar su personalidad calmada. Diminutivos.'

$ws.Range("A41").Value = 'R-7b'
$ws.Range("B41").Value = 'es que cocinar me relaja un monton 😍'
$ws.Range("C41").Value = 'Enviar inmediatamente despues de R-7a, sin esperar respuesta.'

$ws.Range("A42").Value = 'R-7a'
$ws.Range("B42").Value = 'yo estuve haciendo ejercicio un ratico y despues me puse a cocinar algo rico... tururu ✨'
$ws.Range("C42").Value = 'Tururu moment. Diminutivo ''ratico''. Cotidianidad.'

$ws.Range("A43").Value = 'R-6'
$ws.Range("B43").Value = 'oye amor y que hiciste hoy? cuentame despacito que quiero escucharte 💗'
$ws.Range("C43").Value = 'Preguntar por su dia. Diminutivo ''despacito''.'

$ws.Range("A44").Value = 'R-5'
$ws.Range("B44").Value = 'ay me encanta eso 🥰 se nota que eres de buen corazon... eso me llama mucho la atencion ✨'
$ws.Range("C44").Value = 'Validar lo que dice. Ego boost suave.'

$ws.Range("A45").Value = 'R-4b'
$ws.Range("B45").Value = 'tu tienes mascotas, mor?'
$ws.Range("C45").Value = 'Enviar inmediatamente despues de R-4a, sin esperar respuesta.'

$ws.Range("A46").Value = 'R-4a'
$ws.Range("B46").Value = 'alla en la finca de mi familia habia caballos y yo me la pasaba con ellos... sueño con tener uno propio algun dia 💗'
$ws.Range("C46").Value = 'Compartir algo personal genuino. Crear conexion.'

$ws.Range("A47").Value = 'R-3b'
$ws.Range("B47").Value = 'fue dificil dejar todo pero tambien fue lo mas valiente que he hecho'
$ws.Range("C47").Value = 'Enviar inmediatamente despues de R-3a, sin esperar respuesta.'

$ws.Range("A48").Value = 'R-3a'
$ws.Range("B48").Value = 'que lindo ✨ yo naci en Venezuela pero me vine a Colombia hace un tiempito...'
$ws.Range("C48").Value = 'Si nombra un lugar que ella conoce, conectar.'

$ws.Range("A49").Value = 'R-2'
$ws.Range("B49").Value = 'ay que bonito 🥰 oye y tu de donde eres? es que quiero saber todito de ti'
$ws.Range("C49").Value = 'Interes genuino. Diminutivo ''todito''. Devolver pregunta.'

$ws.Range("A50").Value = 'R-1'
$ws.Range("B50").Value = 'hola amor 💗 que alegria verte por aqui... como te sientes hoy?'
$ws.Range("C50").Value = 'Agregar su NOMBRE si se conoce. NUNCA decir cariño/cielo/bebe/princesa. SOLO mor/amor/papi.'
